# Auto-generated edit script: updates market-price / profit data cells
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# hunk 0
$ws.Range("H19").Value = 756.8125
$ws.Range("I19").Value = 507.14285
$ws.Range("J19").Value = 951
$ws.Range("K19").Value = 507.14285
$ws.Range("L19").Value = 951
$ws.Range("M19").Value = -332.14285
$ws.Range("N19").Value = -1301
# hunk 1
$ws.Range("H41").Value = 1573.3462
$ws.Range("I41").Value = 1012.1111
$ws.Range("J41").Value = 2836.125
$ws.Range("K41").Value = 1012.1111
$ws.Range("L41").Value = 2836.125
$ws.Range("M41").Value = -572.1111
$ws.Range("N41").Value = -3716.125
# hunk 2
$ws.Range("H57").Value = 61450
$ws.Range("J57").Value = 61450
$ws.Range("L57").Value = 184350
$ws.Range("N57").Value = -185348
# hunk 3
$ws.Range("H116").Value = 7455.3237
$ws.Range("I116").Value = 7399.6665
$ws.Range("K116").Value = 7399.6665
$ws.Range("M116").Value = -3957.6665
# hunk 4
$ws.Range("H132").Value = 1706.0476
$ws.Range("I132").Value = 1375.2
$ws.Range("K132").Value = 4125.6
$ws.Range("M132").Value = -1595.6
# hunk 5
$ws.Range("H137").Value = 2103.4243
$ws.Range("I137").Value = 1503.8948
$ws.Range("J137").Value = 2917.0715
$ws.Range("K137").Value = 4511.6844
$ws.Range("L137").Value = 8751.2145
$ws.Range("M137").Value = -1961.6844
$ws.Range("N137").Value = -13851.2145
# hunk 6
$ws.Range("H138").Value = 3055.6667
$ws.Range("I138").Value = 2281.65
$ws.Range("J138").Value = 3759.318
$ws.Range("K138").Value = 6844.950000000001
$ws.Range("L138").Value = 11277.954
$ws.Range("M138").Value = -1704.950000000001
$ws.Range("N138").Value = -21557.954

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# hunk 7
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
# hunk 8
$ws.Range("H32").Value = 3313.551
$ws.Range("I32").Value = 2775.0466
$ws.Range("K32").Value = 2775.0466
$ws.Range("M32").Value = -2488.0466
# hunk 9
$ws.Range("H74").Value = 12352100
$ws.Range("I74").Value = 13895238
$ws.Range("K74").Value = 13895238
$ws.Range("M74").Value = -13894364
# hunk 10
$ws.Range("H77").Value = 12352100
$ws.Range("I77").Value = 13895238
$ws.Range("K77").Value = 69476190
$ws.Range("M77").Value = -69471822
# hunk 11
$ws.Range("H92").Value = 19543
$ws.Range("J92").Value = 19543
$ws.Range("L92").Value = 19543
$ws.Range("N92").Value = -24535
# hunk 12
$ws.Range("H122").Value = 2380.125
$ws.Range("I122").Value = 1652
$ws.Range("J122").Value = 4888.1113
$ws.Range("K122").Value = 4956
$ws.Range("L122").Value = 14664.3339
$ws.Range("M122").Value = -2506
$ws.Range("N122").Value = -19564.3339
# hunk 13
$ws.Range("H132").Value = 4361.2173
$ws.Range("I132").Value = 2360.8125
$ws.Range("J132").Value = 8933.571
$ws.Range("K132").Value = 7082.4375
$ws.Range("L132").Value = 26800.713
$ws.Range("M132").Value = -4552.4375
$ws.Range("N132").Value = -31860.713

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# hunk 14
$ws.Range("H81").Value = 49877.668
$ws.Range("J81").Value = 49877.668
$ws.Range("L81").Value = 49877.668
$ws.Range("N81").Value = -51999.668
# hunk 15
$ws.Range("H84").Value = 49877.668
$ws.Range("J84").Value = 49877.668
$ws.Range("L84").Value = 149633.004
$ws.Range("N84").Value = -160241.004
# hunk 16
$ws.Range("H105").Value = 22219.076
$ws.Range("I105").Value = 24827.666
$ws.Range("K105").Value = 24827.666
$ws.Range("M105").Value = -23080.666
# hunk 17
$ws.Range("H134").Value = 4711.75
$ws.Range("I134").Value = 1928.4286
$ws.Range("K134").Value = 5785.2858
$ws.Range("M134").Value = -3250.2858

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# hunk 18
$ws.Range("H16").Value = 587.875
$ws.Range("I16").Value = 672.2
$ws.Range("J16").Value = 447.33334
$ws.Range("K16").Value = 672.2
$ws.Range("L16").Value = 447.33334
$ws.Range("M16").Value = -385.2
$ws.Range("N16").Value = -1021.33334
# hunk 19
$ws.Range("H31").Value = 25399.674
$ws.Range("I31").Value = 2404.5
$ws.Range("J31").Value = 61707.844
$ws.Range("K31").Value = 2404.5
$ws.Range("L31").Value = 61707.844
$ws.Range("M31").Value = -2109.5
$ws.Range("N31").Value = -62297.844
# hunk 20
$ws.Range("H34").Value = 25399.674
$ws.Range("I34").Value = 2404.5
$ws.Range("J34").Value = 61707.844
$ws.Range("K34").Value = 2404.5
$ws.Range("L34").Value = 61707.844
$ws.Range("M34").Value = -2202.5
$ws.Range("N34").Value = -62111.844
# hunk 21
$ws.Range("H113").Value = 587.875
$ws.Range("I113").Value = 672.2
$ws.Range("J113").Value = 447.33334
$ws.Range("K113").Value = 672.2
$ws.Range("L113").Value = 447.33334
$ws.Range("M113").Value = 1497.8
$ws.Range("N113").Value = -4787.33334
# hunk 22
$ws.Range("H132").Value = 4640.325
$ws.Range("I132").Value = 4224.7144
$ws.Range("K132").Value = 12674.1432
$ws.Range("M132").Value = -10144.1432
# hunk 23
$ws.Range("H134").Value = 7750.244
$ws.Range("I134").Value = 3797.238
$ws.Range("K134").Value = 11391.714
$ws.Range("M134").Value = -8856.714

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# hunk 24
$ws.Range("H12").Value = 495.46155
$ws.Range("I12").Value = 139.75
$ws.Range("J12").Value = 653.55554
$ws.Range("K12").Value = 419.25
$ws.Range("L12").Value = 1960.66662
$ws.Range("M12").Value = -246.25
$ws.Range("N12").Value = -2306.66662
# hunk 25
$ws.Range("H42").Value = 7528
$ws.Range("I42").Value = 30
$ws.Range("J42").Value = 11277
$ws.Range("K42").Value = 90
$ws.Range("L42").Value = 33831
$ws.Range("M42").Value = 444
$ws.Range("N42").Value = -34899
# hunk 26
$ws.Range("H80").Value = 4901
$ws.Range("I80").Value = 2468
$ws.Range("J80").Value = 7334
$ws.Range("K80").Value = 7404
$ws.Range("L80").Value = 22002
$ws.Range("M80").Value = -6468
$ws.Range("N80").Value = -23874
# hunk 27
$ws.Range("H83").Value = 4901
$ws.Range("I83").Value = 2468
$ws.Range("J83").Value = 7334
$ws.Range("K83").Value = 22212
$ws.Range("L83").Value = 66006
$ws.Range("M83").Value = -17532
$ws.Range("N83").Value = -75366
# hunk 28
$ws.Range("H122").Value = 1497.8334
$ws.Range("J122").Value = 1742.6
$ws.Range("L122").Value = 15683.4
$ws.Range("N122").Value = -20583.4
# hunk 29
$ws.Range("H138").Value = 5918.615
$ws.Range("I138").Value = 3317.3333
$ws.Range("J138").Value = 11771.5
$ws.Range("K138").Value = 9951.999899999999
$ws.Range("L138").Value = 35314.5
$ws.Range("M138").Value = -4811.999899999999
$ws.Range("N138").Value = -45594.5
# hunk 30
$ws.Range("H139").Value = 6806.8
$ws.Range("I139").Value = 2143.5
$ws.Range("J139").Value = 9915.666999999999
$ws.Range("K139").Value = 6430.5
$ws.Range("L139").Value = 29747.001
$ws.Range("M139").Value = -1290.5
$ws.Range("N139").Value = -40027.001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# hunk 31
$ws.Range("H80").Value = 8000.6665
$ws.Range("I80").Value = 5666.6665
$ws.Range("J80").Value = 12668.667
$ws.Range("K80").Value = 5666.6665
$ws.Range("L80").Value = 12668.667
$ws.Range("M80").Value = -4668.6665
$ws.Range("N80").Value = -14664.667
# hunk 32
$ws.Range("H83").Value = 8000.6665
$ws.Range("I83").Value = 5666.6665
$ws.Range("J83").Value = 12668.667
$ws.Range("K83").Value = 28333.3325
$ws.Range("L83").Value = 63343.335
$ws.Range("M83").Value = -23341.3325
$ws.Range("N83").Value = -73327.33499999999
# hunk 33
$ws.Range("H113").Value = 5658.364
$ws.Range("J113").Value = 5824.2
$ws.Range("L113").Value = 5824.2
$ws.Range("N113").Value = -10164.2
# hunk 34
$ws.Range("H132").Value = 84930.234
$ws.Range("I132").Value = 147648.86
$ws.Range("K132").Value = 442946.58
$ws.Range("M132").Value = -440416.58

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# hunk 35
$ws.Range("H61").Value = 10916.333
$ws.Range("I61").Value = 3874.5
$ws.Range("K61").Value = 3874.5
$ws.Range("M61").Value = -3672.5
# hunk 36
$ws.Range("H93").Value = 5771
$ws.Range("I93").Value = 1662
$ws.Range("J93").Value = 9366.375
$ws.Range("K93").Value = 1662
$ws.Range("L93").Value = 9366.375
$ws.Range("M93").Value = -414
$ws.Range("N93").Value = -11862.375
# hunk 37
$ws.Range("H113").Value = 10916.333
$ws.Range("I113").Value = 3874.5
$ws.Range("K113").Value = 3874.5
$ws.Range("M113").Value = -1704.5
# hunk 38
$ws.Range("H132").Value = 3426.5
$ws.Range("I132").Value = 1977.4615
$ws.Range("J132").Value = 7194
$ws.Range("K132").Value = 5932.3845
$ws.Range("L132").Value = 21582
$ws.Range("M132").Value = -3402.3845
$ws.Range("N132").Value = -26642

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# hunk 39
$ws.Range("H81").Value = 2843.0908
$ws.Range("I81").Value = 1837.0667
$ws.Range("J81").Value = 4998.857
$ws.Range("K81").Value = 3674.1334
$ws.Range("L81").Value = 9997.714
$ws.Range("M81").Value = -2613.1334
$ws.Range("N81").Value = -12119.714
# hunk 40
$ws.Range("H84").Value = 2843.0908
$ws.Range("I84").Value = 1837.0667
$ws.Range("J84").Value = 4998.857
$ws.Range("K84").Value = 18370.667
$ws.Range("L84").Value = 49988.57
$ws.Range("M84").Value = -13066.667
$ws.Range("N84").Value = -60596.57
# hunk 41
$ws.Range("H113").Value = 623.75
$ws.Range("I113").Value = 498.33334
$ws.Range("K113").Value = 1495.00002
$ws.Range("M113").Value = 674.9999800000001
# hunk 42
$ws.Range("H119").Value = 58486.5
$ws.Range("J119").Value = 58486.5
$ws.Range("L119").Value = 58486.5
$ws.Range("N119").Value = -68162.5
# hunk 43
$ws.Range("H132").Value = 5602.6665
$ws.Range("I132").Value = 5260.5454
$ws.Range("J132").Value = 6543.5
$ws.Range("K132").Value = 15781.6362
$ws.Range("L132").Value = 19630.5
$ws.Range("M132").Value = -13251.6362
$ws.Range("N132").Value = -24690.5
